$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('J2').Value = 6019
$ws.Range('J3').Value = 6442
$ws.Range('I4').Value = 1432
$ws.Range('J4').Value = 1391
$ws.Range('J5').Value = 493
$ws.Range('J6').Value = 8336
$ws.Range('I7').Value = 20379
$ws.Range('J7').Value = 22681

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('J3').Value = 41
$ws.Range('J6').Value = 209
$ws.Range('J7').Value = 317

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('J2').Value = 386
$ws.Range('J3').Value = 434
$ws.Range('J6').Value = 492
$ws.Range('J7').Value = 1425

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('J3').Value = 348
$ws.Range('J6').Value = 365
$ws.Range('J7').Value = 1046

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('J3').Value = 118
$ws.Range('J6').Value = 74
$ws.Range('J7').Value = 330

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('J3').Value = 238
$ws.Range('J5').Value = 28
$ws.Range('J6').Value = 201
$ws.Range('J7').Value = 695

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('J6').Value = 92
$ws.Range('J7').Value = 353

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('J2').Value = 179
$ws.Range('J5').Value = 71
$ws.Range('J6').Value = 168
$ws.Range('J7').Value = 671
$ws.Range('J8').Value = 1425
$ws.Range('I11').Value = 298
$ws.Range('J14').Value = 115
$ws.Range('J15').Value = 252
$ws.Range('J19').Value = 670
$ws.Range('J20').Value = 476
$ws.Range('J27').Value = 140
$ws.Range('J29').Value = 1251
$ws.Range('J31').Value = 205
$ws.Range('J33').Value = 1046
$ws.Range('J36').Value = 310
$ws.Range('J37').Value = 695
$ws.Range('J41').Value = 150
$ws.Range('J42').Value = 958
$ws.Range('J44').Value = 172
$ws.Range('J48').Value = 269
$ws.Range('J49').Value = 151
$ws.Range('J50').Value = 135
$ws.Range('J51').Value = 286
$ws.Range('J52').Value = 566
$ws.Range('J53').Value = 317
$ws.Range('J54').Value = 438
$ws.Range('J55').Value = 321
$ws.Range('J61').Value = 24
$ws.Range('J63').Value = 85
$ws.Range('J67').Value = 863
$ws.Range('J76').Value = 345
$ws.Range('J78').Value = 276
$ws.Range('J79').Value = 647
$ws.Range('J84').Value = 190
$ws.Range('J85').Value = 932
$ws.Range('J89').Value = 297
$ws.Range('J91').Value = 259
$ws.Range('J95').Value = 330
$ws.Range('J97').Value = 202
$ws.Range('J98').Value = 166
$ws.Range('J99').Value = 353
$ws.Range('I101').Value = 20379
$ws.Range('J101').Value = 22681

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('J3').Value = 55
$ws.Range('J4').Value = 12
$ws.Range('J6').Value = 58
$ws.Range('J7').Value = 205

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('J2').Value = 215
$ws.Range('J3').Value = 326
$ws.Range('J6').Value = 233
$ws.Range('J7').Value = 863

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('J3').Value = 62
$ws.Range('J7').Value = 190

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('J4').Value = 11
$ws.Range('J7').Value = 151

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('J2').Value = 105
$ws.Range('J7').Value = 438

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('J2').Value = 384
$ws.Range('J3').Value = 433
$ws.Range('J5').Value = 48
$ws.Range('J6').Value = 319
$ws.Range('J7').Value = 1251

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('J2').Value = 45
$ws.Range('J3').Value = 48
$ws.Range('J6').Value = 134
$ws.Range('J7').Value = 269

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('J2').Value = 159
$ws.Range('J3').Value = 196
$ws.Range('J7').Value = 670

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('J2').Value = 53
$ws.Range('J7').Value = 172

$ws = $wb.Worksheets.Item('River North')
$ws.Range('J6').Value = 192
$ws.Range('J7').Value = 345

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('J6').Value = 42
$ws.Range('J7').Value = 115

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('J6').Value = 62
$ws.Range('J7').Value = 168

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('J6').Value = 87
$ws.Range('J7').Value = 150

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('J3').Value = 195
$ws.Range('J6').Value = 499
$ws.Range('J7').Value = 958

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('J2').Value = 36
$ws.Range('J3').Value = 31

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('J6').Value = 80
$ws.Range('J7').Value = 276

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('J6').Value = 168
$ws.Range('J7').Value = 321

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('J6').Value = 61
$ws.Range('J7').Value = 259

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('J3').Value = 222
$ws.Range('J4').Value = 38
$ws.Range('J6').Value = 189
$ws.Range('J7').Value = 647

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('J2').Value = 132
$ws.Range('J3').Value = 164
$ws.Range('J7').Value = 476

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('J6').Value = 96
$ws.Range('J7').Value = 310

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('J3').Value = 204
$ws.Range('J7').Value = 671

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('J3').Value = 60
$ws.Range('J7').Value = 252

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('J6').Value = 102
$ws.Range('J7').Value = 166

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('J6').Value = 43
$ws.Range('J7').Value = 135

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('I4').Value = 26
$ws.Range('I7').Value = 298

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('J4').Value = 14
$ws.Range('J7').Value = 179

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('J2').Value = 34
$ws.Range('J7').Value = 202

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('J3').Value = 84
$ws.Range('J7').Value = 297

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range('J2').Value = 24
$ws.Range('J7').Value = 71

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('J3').Value = 33
$ws.Range('J7').Value = 140

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('J6').Value = 116
$ws.Range('J7').Value = 286

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('J3').Value = 333
$ws.Range('J4').Value = 63
$ws.Range('J5').Value = 20
$ws.Range('J6').Value = 271
$ws.Range('J7').Value = 932

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('J3').Value = 172
$ws.Range('J7').Value = 566

$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Range('J2').Value = 12
$ws.Range('J6').Value = 24
